# Apply updated cryptocurrency price/volume figures to Sheet1 (2024-09-21 GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell D4 keeps its untouched, default (style 0) formatting; we use it to restore
# the plain format on D-column cells after forcing a Text entry below, so that numeric-
# looking price strings (e.g. "586.37") are not auto-converted to the Number type.
$ws.Range("D4").Copy() | Out-Null

$ws.Range("D2").Value = "63.096.79"
$ws.Range("E2").Value = "  -0.62%  "

$ws.Range("D3").Value = "2.552.21"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.37"
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = "  +2.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.54"
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").Value = "  -2.48%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -0.49%  "

$ws.Range("E9").Value = "  -0.53%  "

$ws.Range("E10").Value = "  -3.44%  "

$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("E12").Value = "  -1.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.55"
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("E13").Value = "  -3.01%  "

$ws.Range("D14").Value = "3.005.13"
$ws.Range("E14").Value = "  +0.13%  "

$ws.Range("D15").Value = "63.005.53"
$ws.Range("E15").Value = "  -0.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000143"
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").Value = "2.559.42"
$ws.Range("E17").Value = "  +2.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.38"
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "  -2.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "337.25"
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").Value = "  -1.02%  "

$ws.Range("E20").Value = "  -0.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.77"
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("E21").Value = "  -1.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.75"
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = "  -0.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.169"
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("E24").Value = "  -0.46%  "

$ws.Range("E25").Value = "  +1.02%  "

$ws.Range("E26").Value = "  -1.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("E29").Value = "  +8.94%  "

$ws.Range("E30").Value = "  +6.35%  "

$ws.Range("D31").Value = "0.0₃0817"
$ws.Range("E31").Value = "  -2.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "178.71"
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("E32").Value = "  +0.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "421.37"
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33").Value = "  +0.51%  "

$ws.Range("E34").Value = "  -0.74%  "

$ws.Range("E35").Value = "  -1.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.15"
$ws.Range("D36").PasteSpecial(-4122)

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.38"
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("E38").Value = "  -2.49%  "

$ws.Range("E39").Value = "  -1.08%  "

$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "39.74"
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("E41").Value = "  -0.63%  "

$ws.Range("E42").Value = "  -2.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.79"
$ws.Range("D43").PasteSpecial(-4122)
$ws.Range("E43").Value = "  -0.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.88"
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("E44").Value = "  -1.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0543"
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("E45").Value = "  +2.17%  "

$ws.Range("E46").Value = "  -1.40%  "

$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("E48").Value = "  -1.78%  "

$ws.Range("E49").Value = "  -1.68%  "

$ws.Range("E50").Value = "  -5.86%  "

$ws.Range("E51").Value = "  -0.39%  "

$excel.CutCopyMode = $false
